# The deck's theme (ppt/theme/theme1.xml, used by the slide master / all
# slides) switches from the custom "Integral" / Red Violet palette to the
# stock PowerPoint "Office Theme" palette (the palette that used to sit,
# unused by any slide, on the secondary theme part used only by the notes
# master). Re-apply each of the 12 theme colour-scheme slots on the slide
# master's theme to the Office Theme RGB values.

$p = $ppt.ActivePresentation

# Office Theme colours (RRGGBB), in clrScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

function ConvertTo-BgrInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$scheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $bgr = ConvertTo-BgrInt $officeColors[$i - 1]
    $scheme.Colors($i).RGB = $bgr
}

Write-Output "Theme palette updated to Office Theme"
